$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'womans tights winter'
$ws.Cells.Item(2, 1).Value = 'womans under armour'
$ws.Cells.Item(3, 1).Value = 'womans winter running gear'
$ws.Cells.Item(4, 1).Value = 'womans workout clothing'
$ws.Cells.Item(5, 1).Value = 'womans workout gear'
$ws.Cells.Item(6, 1).Value = 'wome nike'
$ws.Cells.Item(7, 1).Value = 'womems leggings'
$ws.Cells.Item(8, 1).Value = 'womems tights'
$ws.Cells.Item(9, 1).Value = 'women 7 inch shorts'
$ws.Cells.Item(10, 1).Value = 'women active wear pants'
$ws.Cells.Item(11, 1).Value = 'women all weather jacket'
$ws.Cells.Item(12, 1).Value = 'women all weather jackets'
$ws.Cells.Item(13, 1).Value = 'women athleta'
$ws.Cells.Item(14, 1).Value = 'women athletic capri'
$ws.Cells.Item(15, 1).Value = 'women athletic capris'
$ws.Cells.Item(16, 1).Value = 'women athletic clothes'
$ws.Cells.Item(17, 1).Value = 'women athletic gear'
$ws.Cells.Item(18, 1).Value = 'women athletic leggings'
$ws.Cells.Item(19, 1).Value = 'women athletic leggings capri'
$ws.Cells.Item(20, 1).Value = 'women athletic tights'
$ws.Cells.Item(21, 1).Value = 'women athletic wear'
$ws.Cells.Item(22, 1).Value = 'women basketball'
$ws.Cells.Item(23, 1).Value = 'women basketball pants'
$ws.Cells.Item(24, 1).Value = 'women basketball shorts'
$ws.Cells.Item(25, 1).Value = 'women best leggings'
$ws.Cells.Item(26, 1).Value = 'women black compression shorts'
$ws.Cells.Item(27, 1).Value = 'women black leggings'
$ws.Cells.Item(28, 1).Value = 'women black pants'
$ws.Cells.Item(29, 1).Value = 'women black stockings'
$ws.Cells.Item(30, 1).Value = 'women black tights'
$ws.Cells.Item(31, 1).Value = 'women boxing'
$ws.Cells.Item(32, 1).Value = 'women boxing gear'
$ws.Cells.Item(33, 1).Value = 'women boxing shorts'
$ws.Cells.Item(34, 1).Value = 'women breathable leggings'
$ws.Cells.Item(35, 1).Value = 'women brown tights'
$ws.Cells.Item(36, 1).Value = 'women buttery soft leggings'
$ws.Cells.Item(37, 1).Value = 'women capri'
$ws.Cells.Item(38, 1).Value = 'women capri leggings'
$ws.Cells.Item(39, 1).Value = 'women capri overalls'
$ws.Cells.Item(40, 1).Value = 'women capri pants'
$ws.Cells.Item(41, 1).Value = 'women capri shorts'
$ws.Cells.Item(42, 1).Value = 'women capri workout leggings'
$ws.Cells.Item(43, 1).Value = 'women capri workout pants'
$ws.Cells.Item(44, 1).Value = 'women capri yoga pants'
$ws.Cells.Item(45, 1).Value = 'women capris'
$ws.Cells.Item(46, 1).Value = 'women capris leggings'
$ws.Cells.Item(47, 1).Value = 'women capris pants'
$ws.Cells.Item(48, 1).Value = 'women club dresses 2017'
$ws.Cells.Item(49, 1).Value = 'women cold gear'
$ws.Cells.Item(50, 1).Value = 'women cold gear running'
$ws.Cells.Item(51, 1).Value = 'women cold weather running gear'
$ws.Cells.Item(52, 1).Value = 'women compression'
$ws.Cells.Item(53, 1).Value = 'women compression brace'
$ws.Cells.Item(54, 1).Value = 'women compression capri'
$ws.Cells.Item(55, 1).Value = 'women compression capri leggings'
$ws.Cells.Item(56, 1).Value = 'women compression clothing'
$ws.Cells.Item(57, 1).Value = 'women compression gear'
$ws.Cells.Item(58, 1).Value = 'women compression hose'
$ws.Cells.Item(59, 1).Value = 'women compression jacket'
$ws.Cells.Item(60, 1).Value = 'women compression knee highs'
$ws.Cells.Item(61, 1).Value = 'women compression leggings'
$ws.Cells.Item(62, 1).Value = 'women compression leggings tall'
$ws.Cells.Item(63, 1).Value = 'women compression panties'
$ws.Cells.Item(64, 1).Value = 'women compression pants cold gear'
$ws.Cells.Item(65, 1).Value = 'women compression pantyhose'
$ws.Cells.Item(66, 1).Value = 'women compression running pants'
$ws.Cells.Item(67, 1).Value = 'women compression running tight'
$ws.Cells.Item(68, 1).Value = 'women compression shorts'
$ws.Cells.Item(69, 1).Value = 'women compression shorts pack'
$ws.Cells.Item(70, 1).Value = 'women compression shorts with pocket'
$ws.Cells.Item(71, 1).Value = 'women compression tank'
$ws.Cells.Item(72, 1).Value = 'women compression tight'
$ws.Cells.Item(73, 1).Value = 'women compression tights'
$ws.Cells.Item(74, 1).Value = 'women compression tights for running'
$ws.Cells.Item(75, 1).Value = 'women compression top'
$ws.Cells.Item(76, 1).Value = 'women compression underwear'
$ws.Cells.Item(77, 1).Value = 'women compression waist'
$ws.Cells.Item(78, 1).Value = 'women compression wear'
$ws.Cells.Item(79, 1).Value = 'women cwx'
$ws.Cells.Item(80, 1).Value = 'women cycle shorts'
$ws.Cells.Item(81, 1).Value = 'women cycling pants'
$ws.Cells.Item(82, 1).Value = 'women cycling shorts'
$ws.Cells.Item(83, 1).Value = 'women down pants'
$ws.Cells.Item(84, 1).Value = 'women fashion tights'
$ws.Cells.Item(85, 1).Value = 'women fatigue pants'
$ws.Cells.Item(86, 1).Value = 'women fitness clothes'
$ws.Cells.Item(87, 1).Value = 'women fitness underwear'
$ws.Cells.Item(88, 1).Value = 'women fitted ski pants'
$ws.Cells.Item(89, 1).Value = 'women gym leggings'
$ws.Cells.Item(90, 1).Value = 'women gym panties'
$ws.Cells.Item(91, 1).Value = 'women gym pants'
$ws.Cells.Item(92, 1).Value = 'women gym short'
$ws.Cells.Item(93, 1).Value = 'women gym tights'
$ws.Cells.Item(94, 1).Value = 'women high top nike'
$ws.Cells.Item(95, 1).Value = 'women high waist'
$ws.Cells.Item(96, 1).Value = 'women high waist pants'
$ws.Cells.Item(97, 1).Value = 'women hiking capris'
$ws.Cells.Item(98, 1).Value = 'women hiking tights'
$ws.Cells.Item(99, 1).Value = 'women insulated pants'
$ws.Cells.Item(100, 1).Value = 'women jacket all weather'
